$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new requirement row (rml-io-r7) right after the existing
# rml-io-r6 entry, mirroring the Identifier / Requirement / Provenance
# columns used by the rest of the table.
$ws.Range("A8").Value = "rml-io-r7"
$ws.Range("B8").Value = "Logical sources and logical targets may indicate relative paths to resources"
$ws.Range("C8").Value = "rml-io"

# The Requirement column wraps text (style already applied to B8), so the
# row needs the same taller row height used by the other requirement rows.
$ws.Rows.Item(8).RowHeight = 34

# The author's last selection in the sheet moved down to B9 after filling
# in row 8.
$ws.Range("B9").Select()
